# Add the new vocabulary words (List8 / List9 groups) to Sheet1, column A,
# continuing directly after the existing data (through row 369).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$words = @(
    "benign",
    "malignant",
    "tumor",
    "malefactor",
    "benediction",
    "malediction",
    "benevolent",
    "malevolent",
    "malice",
    "malign",
    "malicious",
    "malfeasance",
    "malcontent",
    "melody",
    "malady",
    "malaise",
    "List8",
    "belie",
    "belittle",
    "beguile",
    "guile",
    "guileless",
    "besiege",
    "beleaguer",
    "league",
    "bestow",
    "endow",
    "donate",
    "dote",
    "condone",
    "profane",
    "blaspheme",
    "blantant",
    "bloated",
    "blotch",
    "boycott",
    "flamboyant",
    "boisterous",
    "lobster",
    "bolster",
    "bind",
    "bondage",
    "bandage",
    "enthrall",
    "thrall",
    "slave",
    "bondman",
    "enfranchise",
    "bravado",
    "wig",
    "denture",
    "brave",
    "courageous",
    "gallant",
    "intripid",
    "tripid",
    "tremble",
    "plucky",
    "luck",
    "breach",
    "breed",
    "brood",
    "embryo",
    "braid",
    "broach",
    "List9"
)

# Row 1 already uses the bold red "list header" style (s=3 in the original
# sheet); reuse its formatting for the new "List8"/"List9" header rows so we
# do not introduce any new font/style table entries.
$headerTemplate = $ws.Cells.Item(1, 1)

$startRow = 370
for ($i = 0; $i -lt $words.Count; $i++) {
    $r = $startRow + $i
    $word = $words[$i]
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = $word
    if ($word -like "List*") {
        $headerTemplate.Copy() | Out-Null
        $cell.PasteSpecial(-4122) | Out-Null
    }
}
$excel.CutCopyMode = $false

# Move the view/selection to reflect where the user ended up after the edit.
$excel.ActiveWindow.ScrollRow = 423
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C436").Select() | Out-Null
